$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("out_vars")

$row = 76
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2020-08-14"
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 2).Value = 511369
$ws.Cells.Item($row, 3).Value = 559974
$ws.Cells.Item($row, 4).Value = 85509
$ws.Cells.Item($row, 5).Value = 55908
$ws.Cells.Item($row, 6).Value = 26.33
